$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill previously-empty I2:L2 with computed numeric values
$ws.Range("I2").Value = -0.004430112316335983
$ws.Range("J2").Value = 0.2432629225634229
$ws.Range("K2").Value = -0.04778595182768761
$ws.Range("L2").Value = 2.671675381450449

# Row 20: fill previously-empty I20:L20 with computed numeric values
$ws.Range("I20").Value = -0.1807720062054129
$ws.Range("J20").Value = 0.3297608760644796
$ws.Range("K20").Value = 0.05159556685682618
$ws.Range("L20").Value = 2.232560601757156
